$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 281; existing rows 281..355 shift down to 282..356.
$ws.Rows.Item(281).Insert()

# Populate the newly inserted row 281 with its data.
$ws.Range("A281").Value = 4
$ws.Range("B281").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C281").Value = "Los Lagos"
$ws.Range("D281").Value = 44642
$ws.Range("E281").Value = 10
$ws.Range("F281").Value = 100114001
$ws.Range("G281").Value = "Papa"
$ws.Range("H281").Value = "Patagonia"
$ws.Range("I281").Value = "1a (cosecha)"
$ws.Range("J281").Value = 600
$ws.Range("K281").Value = 7000
$ws.Range("L281").Value = 8000
$ws.Range("M281").Value = 7500
$ws.Range("N281").Value = "`$/saco 25 kilos"
$ws.Range("O281").Value = "Provincia de Llanquihue"
$ws.Range("P281").Value = 300
$ws.Range("Q281").Value = 25
$ws.Range("R281").Value = "Hortaliza"
